$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed crypto price / volume(1h) data as plain text,
# matching the source formatting (dot-separated thousands, percentages
# padded with two leading/trailing spaces). A handful of price values
# parse as plain numbers, so force a text number format first to stop
# Excel from silently converting them (e.g. "9.50" -> 9.5).

$ws.Range("D2").Value = "41.084.69"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.425.72"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.59"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.32"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  -2.64%  "
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0836"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.02"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "2.799.87"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.63"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "2.430.81"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "41.016.44"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.46"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.89"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.16"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.57"
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.25"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.95"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.88"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("E39").Value = "  -3.96%  "
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("E42").Value = "  -5.07%  "
$ws.Range("D43").Value = "1.995.75"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.55"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.50"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "2.662.72"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.62"
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.30"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.97"
$ws.Range("E51").Value = "  -0.69%  "
